$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BassFly_uHat.bom")

# Fix 0.1mm misalignment of mounting holes:
# Update MOUSER (H) and DIGIKEY (I) columns for row 9 (J4 - OLED module)
# from "N.M." to "-"
$ws.Range("H9").Value = "-"
$ws.Range("I9").Value = "-"

# Reflect the final selected cell as left by the author when saving
$ws.Range("H28").Select()

$wb.Save()
